# Removed an obsolete test suite ("JobStatus") from the Sheet1 summary
# table. The rows below it (Read, Runtimes) shift up one position and the
# now-unused trailing row is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 ("JobStatus") is replaced by what used to be row 5 ("Read").
$ws.Cells.Item(4, 1).Value = "Read"
$ws.Cells.Item(4, 2).Value = 12
$ws.Cells.Item(4, 3).Value = 12
$ws.Cells.Item(4, 4).Value = "Finished"
$ws.Cells.Item(4, 5).Value = "Contains four partially automated test cases."

# Row 5 is replaced by what used to be row 6 ("Runtimes").
$ws.Cells.Item(5, 1).Value = "Runtimes"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 13
$ws.Cells.Item(5, 4).Value = "Suited to Manual"
$ws.Cells.Item(5, 5).Value = "Timing"

# Row 6 no longer holds a data entry (table shrank by one row); clear it.
$ws.Range("A6:E6").Clear()

# Selection moves to the newly-edited row.
$ws.Range("A4:E4").Select() | Out-Null

# The conditional-formatting range over column D shrinks along with the
# table (D2:D52 -> D2:D51).
$fc = $ws.Range("D2").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D2:D51")) | Out-Null
